$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "Prerequisites"
$ws.Range("D1").Value = "Corequisites"
$ws.Range("E1").Value = "Concurrent"
$ws.Range("F1").Value = "Recommended"
$ws.Range("G1").Value = "Terms Typically Offered"
$ws.Range("C2").Value = "NA"
$ws.Range("D2").Value = "NA"
$ws.Range("E2").Value = "NA"
$ws.Range("F2").Value = "NA"
$ws.Range("G2").Value = "TBD"
$ws.Range("C3").Value = "ENGL 129."
$ws.Range("D3").Value = "NA"
$ws.Range("E3").Value = "ENGL 135."
$ws.Range("F3").Value = "NA"
$ws.Range("G3").Value = "TBD "
$ws.Range("C4").Value = "NA"
$ws.Range("D4").Value = "NA"
$ws.Range("E4").Value = "NA"
$ws.Range("F4").Value = "NA"
$ws.Range("G4").Value = "TBD"
$ws.Range("C5").Value = "ENGL 131."
$ws.Range("D5").Value = "NA"
$ws.Range("E5").Value = "ENGL 135."
$ws.Range("F5").Value = "NA"
$ws.Range("G5").Value = "TBD "
$ws.Range("C6").Value = "GE A1 eligibility for Written Communication Placement upon admissions."
$ws.Range("D6").Value = "NA"
$ws.Range("E6").Value = "NA"
$ws.Range("F6").Value = "NA"
$ws.Range("G6").Value = "F, W"
$ws.Range("C7").Value = "GE A1 eligibility for Written Communication placement upon admissions."
$ws.Range("D7").Value = "NA"
$ws.Range("E7").Value = "NA"
$ws.Range("F7").Value = "NA"
$ws.Range("G7").Value = "F, W, SP"
$ws.Range("C8").Value = "NA"
$ws.Range("D8").Value = "ENGL 129, ENGL 130, ENGL 131, ENGL 132, ENGL 133, or ENGL 134."
$ws.Range("E8").Value = "NA"
$ws.Range("F8").Value = "NA"
$ws.Range("G8").Value = "F,W,SP,SU"
$ws.Range("C9").Value = "Completion of GE Area A1 with a grade of C- or better."
$ws.Range("D9").Value = "NA"
$ws.Range("E9").Value = "NA"
$ws.Range("F9").Value = "Completion of GE Area A2."
$ws.Range("G9").Value = "F,W,SP,SU "
$ws.Range("C10").Value = "Completion of GE Area A1 with a grade of C- or better."
$ws.Range("D10").Value = "NA"
$ws.Range("E10").Value = "NA"
$ws.Range("F10").Value = "Completion of GE Area A2."
$ws.Range("G10").Value = "F, W, SP "
$ws.Range("C11").Value = "Completion of GE Area A1 with a grade of C- or better. For Engineering students only."
$ws.Range("D11").Value = "NA"
$ws.Range("E11").Value = "NA"
$ws.Range("F11").Value = "Completion of GE Area A2."
$ws.Range("G11").Value = "F,W,SP,SU "
$ws.Range("C12").Value = "Completion of GE Area A1 with a grade of C- or better and consent of instructor."
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("F12").Value = "NA"
$ws.Range("G12").Value = "F,W,SP,SU"
$ws.Range("C13").Value = "Completion of GE A1 with a grade of C- or better; for English majors only."
$ws.Range("D13").Value = "NA"
$ws.Range("E13").Value = "NA"
$ws.Range("F13").Value = "NA"
$ws.Range("G13").Value = "W, SP"
$ws.Range("C14").Value = "Completion of GE Area A with grades of C- or better; and ENGL 251. For English majors only."
$ws.Range("D14").Value = "NA"
$ws.Range("E14").Value = "NA"
$ws.Range("F14").Value = "NA"
$ws.Range("G14").Value = "F, W, SP"
$ws.Range("C15").Value = "ENGL 203; for English majors only."
$ws.Range("D15").Value = "NA"
$ws.Range("E15").Value = "NA"
$ws.Range("F15").Value = "NA"
$ws.Range("G15").Value = "W, SP"
$ws.Range("C16").Value = "ENGL 204; for English majors only."
$ws.Range("D16").Value = "NA"
$ws.Range("E16").Value = "NA"
$ws.Range("F16").Value = "NA"
$ws.Range("G16").Value = "F, SP"
$ws.Range("C17").Value = "Completion of GE Area A with grades of C- or better."
$ws.Range("D17").Value = "NA"
$ws.Range("E17").Value = "NA"
$ws.Range("F17").Value = "NA"
$ws.Range("G17").Value = "TBD"
$ws.Range("C18").Value = "Completion of GE Areas A1 and A3 with a grade of C- or better; and one course in GE Area B."
$ws.Range("D18").Value = "NA"
$ws.Range("E18").Value = "NA"
$ws.Range("F18").Value = "NA"
$ws.Range("G18").Value = "F"
$ws.Range("C19").Value = "Completion of GE Area A with grades of C- or better."
$ws.Range("D19").Value = "NA"
$ws.Range("E19").Value = "NA"
$ws.Range("F19").Value = "NA"
$ws.Range("G19").Value = "TBD"
$ws.Range("C20").Value = "Completion of GE Area A with grades of C- or better."
$ws.Range("D20").Value = "NA"
$ws.Range("E20").Value = "NA"
$ws.Range("F20").Value = "NA"
$ws.Range("G20").Value = "W, SP"
$ws.Range("C21").Value = "Completion of GE Area A with grades of C- or better."
$ws.Range("D21").Value = "NA"
$ws.Range("E21").Value = "NA"
$ws.Range("F21").Value = "NA"
$ws.Range("G21").Value = "F"
$ws.Range("C22").Value = "Completion of GE Area A with grades of C- or better."
$ws.Range("D22").Value = "NA"
$ws.Range("E22").Value = "NA"
$ws.Range("F22").Value = "NA"
$ws.Range("G22").Value = "F,W,SP,SU"
$ws.Range("C23").Value = "Completion of GE Area A with grades of C- or better."
$ws.Range("D23").Value = "NA"
$ws.Range("E23").Value = "NA"
$ws.Range("F23").Value = "NA"
$ws.Range("G23").Value = "F,W,SP,SU"
$ws.Range("C24").Value = "Completion of GE Area A with grades of C- or better."
$ws.Range("D24").Value = "NA"
$ws.Range("E24").Value = "NA"
$ws.Range("F24").Value = "NA"
$ws.Range("G24").Value = "TBD"
$ws.Range("C25").Value = "Completion of GE Area A with grades of C- or better."
$ws.Range("D25").Value = "NA"
$ws.Range("E25").Value = "NA"
$ws.Range("F25").Value = "NA"
$ws.Range("G25").Value = "TBD"
$ws.Range("C26").Value = "Completion of GE Area A with grades of C- or better."
$ws.Range("D26").Value = "NA"
$ws.Range("E26").Value = "NA"
$ws.Range("F26").Value = "NA"
$ws.Range("G26").Value = "TBD"
$ws.Range("C27").Value = "Open to undergraduate students and consent of instructor."
$ws.Range("D27").Value = "NA"
$ws.Range("E27").Value = "NA"
$ws.Range("F27").Value = "NA"
$ws.Range("G27").Value = "TBD"
$ws.Range("C28").Value = "Completion of GE Area A with grades of C- or better."
$ws.Range("D28").Value = "NA"
$ws.Range("E28").Value = "NA"
$ws.Range("F28").Value = "NA"
$ws.Range("G28").Value = "F, W"
$ws.Range("C29").Value = "ENGL 205; for English majors only."
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = "NA"
$ws.Range("F29").Value = "NA"
$ws.Range("G29").Value = "F, W"
$ws.Range("C30").Value = "ENGL 303; for English majors only."
$ws.Range("D30").Value = "NA"
$ws.Range("E30").Value = "NA"
$ws.Range("F30").Value = "NA"
$ws.Range("G30").Value = "W, SP"
$ws.Range("C31").Value = "ENGL 304; for English majors only."
$ws.Range("D31").Value = "NA"
$ws.Range("E31").Value = "NA"
$ws.Range("F31").Value = "NA"
$ws.Range("G31").Value = "F, SP"
$ws.Range("C32").Value = "ENGL 305; for English majors only."
$ws.Range("D32").Value = "NA"
$ws.Range("E32").Value = "NA"
$ws.Range("F32").Value = "NA"
$ws.Range("G32").Value = "F, SP"
$ws.Range("C33").Value = "Completion of GE Area A with grades of C- or better."
$ws.Range("D33").Value = "NA"
$ws.Range("E33").Value = "NA"
$ws.Range("F33").Value = "Junior standing."
$ws.Range("G33").Value = "F,W,SP,SU "
$ws.Range("C34").Value = "Junior standing; completion of GE Area A with grades of C- or better; completion of one course in GE Area B1 with a grade of C- or better; and completion of GE Area C1."
$ws.Range("D34").Value = "NA"
$ws.Range("E34").Value = "NA"
$ws.Range("F34").Value = "NA"
$ws.Range("G34").Value = "TBD"
$ws.Range("C35").Value = "Junior standing; completion of GE Area A with grades of C- or better; completion of one course in GE Area B1 with a grade of C- or better; and completion of GE Area C1."
$ws.Range("D35").Value = "NA"
$ws.Range("E35").Value = "NA"
$ws.Range("F35").Value = "ENGL 133."
$ws.Range("G35").Value = "W "
$ws.Range("C36").Value = "Completion of GE Area A with grades of C- or better."
$ws.Range("D36").Value = "NA"
$ws.Range("E36").Value = "NA"
$ws.Range("F36").Value = "Junior standing."
$ws.Range("G36").Value = "W "
$ws.Range("C37").Value = "Completion of GE Area A with grades of C- or better."
$ws.Range("D37").Value = "NA"
$ws.Range("E37").Value = "NA"
$ws.Range("F37").Value = "Junior standing."
$ws.Range("G37").Value = "SP "
$ws.Range("C38").Value = "Junior standing or English major; completion of GE Area A with grades of C- or better; completion of one course in GE Area B1 with a grade of C- or better; and completion of GE Area C1."
$ws.Range("D38").Value = "NA"
$ws.Range("E38").Value = "NA"
$ws.Range("F38").Value = "NA"
$ws.Range("G38").Value = "TBD"
$ws.Range("C39").Value = "Junior standing or English major; completion of GE Area A with grades of C- or better; completion of one course in GE Area B1 with a grade of C- or better; and completion of GE Area C1."
$ws.Range("D39").Value = "NA"
$ws.Range("E39").Value = "NA"
$ws.Range("F39").Value = "NA"
$ws.Range("G39").Value = "TBD"
$ws.Range("C40").Value = "Junior standing or English major; completion of GE Area A with grades of C- or better; completion of one course in GE Area B1 with a grade of C- or better; and completion of GE Area C1."
$ws.Range("D40").Value = "NA"
$ws.Range("E40").Value = "NA"
$ws.Range("F40").Value = "NA"
$ws.Range("G40").Value = "TBD"
$ws.Range("C41").Value = "Junior standing or English major; completion of GE Area A with grades of C- or better; completion of one course in GE Area B1 with a grade of C- or better; and completion of GE Area C1."
$ws.Range("D41").Value = "NA"
$ws.Range("E41").Value = "NA"
$ws.Range("F41").Value = "NA"
$ws.Range("G41").Value = "TBD"
$ws.Range("C42").Value = "Junior standing or English major; completion of GE Area A with grades of C- or better; completion of one course in GE Area B1 with a grade of C- or better; and completion of GE Area C1."
$ws.Range("D42").Value = "NA"
$ws.Range("E42").Value = "NA"
$ws.Range("F42").Value = "NA"
$ws.Range("G42").Value = "TBD"
$ws.Range("C43").Value = "Junior standing or English major; completion of GE Area A with grades of C- or better; completion of one course in GE Area B1 with a grade of C- or better; and completion of GE Area C1."
$ws.Range("D43").Value = "NA"
$ws.Range("E43").Value = "NA"
$ws.Range("F43").Value = "NA"
$ws.Range("G43").Value = "F, W"
$ws.Range("C44").Value = "Junior standing or English major; completion of GE Area A with grades of C- or better; completion of one course in GE Area B1 with a grade of C- or better; and completion of GE Area C1."
$ws.Range("D44").Value = "NA"
$ws.Range("E44").Value = "NA"
$ws.Range("F44").Value = "NA"
$ws.Range("G44").Value = "W"
$ws.Range("C45").Value = "Junior standing or English major; completion of GE Area A with grades of C- or better; completion of one course in GE Area B1 with a grade of C- or better; and completion of GE Area C1."
$ws.Range("D45").Value = "NA"
$ws.Range("E45").Value = "NA"
$ws.Range("F45").Value = "NA"
$ws.Range("G45").Value = "TBD"
$ws.Range("C46").Value = "Junior standing or English major; completion of GE Area A with grades of C- or better; completion of one course in GE Area B1 with a grade of C- or better; and completion of GE Area C1."
$ws.Range("D46").Value = "NA"
$ws.Range("E46").Value = "NA"
$ws.Range("F46").Value = "NA"
$ws.Range("G46").Value = "TBD"
$ws.Range("C47").Value = "Junior standing or English major; completion of GE Area A with grades of C- or better; completion of one course in GE Area B1 with a grade of C- or better; and completion of GE Area C1."
$ws.Range("D47").Value = "NA"
$ws.Range("E47").Value = "NA"
$ws.Range("F47").Value = "NA"
$ws.Range("G47").Value = "TBD"
$ws.Range("C48").Value = "Junior standing or English major; completion of GE Area A with grades of C- or better; completion of one course in GE Area B1 with a grade of C- or better; and completion of GE Area C1."
$ws.Range("D48").Value = "NA"
$ws.Range("E48").Value = "NA"
$ws.Range("F48").Value = "NA"
$ws.Range("G48").Value = "TBD"
$ws.Range("C49").Value = "Junior standing or English major; completion of GE Area A with grades of C- or better; completion of one course in GE Area B1 with a grade of C- or better; and completion of GE Area C1."
$ws.Range("D49").Value = "NA"
$ws.Range("E49").Value = "NA"
$ws.Range("F49").Value = "NA"
$ws.Range("G49").Value = "SU"
$ws.Range("C50").Value = "Junior standing or English major; completion of GE Area A with grades of C- or better; completion of one course in GE Area B1 with a grade of C- or better; and completion of GE Area C1."
$ws.Range("D50").Value = "NA"
$ws.Range("E50").Value = "NA"
$ws.Range("F50").Value = "NA"
$ws.Range("G50").Value = "F, W, SP"
$ws.Range("C51").Value = "Junior standing or English major; completion of GE Area A with grades of C- or better; completion of one course in GE Area B1 with a grade of C- or better; and completion of GE Area C1."
$ws.Range("D51").Value = "NA"
$ws.Range("E51").Value = "NA"
$ws.Range("F51").Value = "NA"
$ws.Range("G51").Value = "F"
$ws.Range("C52").Value = "Junior standing or English major; completion of GE Area A with grades of C- or better; completion of one course in GE Area B1 with a grade of C- or better; and completion of GE Area C1."
$ws.Range("D52").Value = "NA"
$ws.Range("E52").Value = "NA"
$ws.Range("F52").Value = "NA"
$ws.Range("G52").Value = "W"
$ws.Range("C53").Value = "Junior standing or English major; completion of GE Area A with grades of C- or better; completion of one course in GE Area B1 with a grade of C- or better; and completion of GE Area C1."
$ws.Range("D53").Value = "NA"
$ws.Range("E53").Value = "NA"
$ws.Range("F53").Value = "NA"
$ws.Range("G53").Value = "TBD"
$ws.Range("C54").Value = "Junior standing or English major; completion of GE Area A with grades of C- or better; completion of one course in GE Area B1 with a grade of C- or better; and completion of GE Area C1."
$ws.Range("D54").Value = "NA"
$ws.Range("E54").Value = "NA"
$ws.Range("F54").Value = "NA"
$ws.Range("G54").Value = "TBD"
$ws.Range("C55").Value = "Junior standing or English major; completion of GE Area A with grades of C- or better; completion of one course in GE Area B1 with a grade of C- or better; and completion of GE Area C1."
$ws.Range("D55").Value = "NA"
$ws.Range("E55").Value = "NA"
$ws.Range("F55").Value = "NA"
$ws.Range("G55").Value = "TBD"
$ws.Range("C56").Value = "Junior standing or English major; completion of GE Area A with grades of C- or better; completion of one course in GE Area B1 with a grade of C- or better; and completion of GE Area C1."
$ws.Range("D56").Value = "NA"
$ws.Range("E56").Value = "NA"
$ws.Range("F56").Value = "NA"
$ws.Range("G56").Value = "TBD"
$ws.Range("C57").Value = "Junior standing or English major; completion of GE Area A with grades of C- or better; completion of one course in GE Area B1 with a grade of C- or better; and completion of GE Area C1."
$ws.Range("D57").Value = "NA"
$ws.Range("E57").Value = "NA"
$ws.Range("F57").Value = "NA"
$ws.Range("G57").Value = "TBD"
$ws.Range("C58").Value = "Junior standing; completion of GE Area A1 with grades of C- or better; and completion of GE Area C1."
$ws.Range("D58").Value = "NA"
$ws.Range("E58").Value = "NA"
$ws.Range("F58").Value = "EDUC 300."
$ws.Range("G58").Value = "W "
$ws.Range("C59").Value = "Junior standing; completion of GE C1 with a grade of ``B' or better, or consent of instructor."
$ws.Range("D59").Value = "NA"
$ws.Range("E59").Value = "NA"
$ws.Range("F59").Value = "NA"
$ws.Range("G59").Value = "F, W"
$ws.Range("C60").Value = "Junior standing or English major; completion of GE Area A with grades of C- or better; completion of one course in GE Area B1 with a grade of C- or better; and completion of GE Area C1."
$ws.Range("D60").Value = "NA"
$ws.Range("E60").Value = "NA"
$ws.Range("F60").Value = "NA"
$ws.Range("G60").Value = "W"
$ws.Range("C61").Value = "Junior standing or English major; completion of GE Area A with grades of C- or better; completion of one course in GE Area B1 with a grade of C- or better; and completion of GE Area C1."
$ws.Range("D61").Value = "NA"
$ws.Range("E61").Value = "NA"
$ws.Range("F61").Value = "NA"
$ws.Range("G61").Value = "TBD"
$ws.Range("C62").Value = "Junior standing or English major; completion of GE Area A with grades of C- or better; completion of one course in GE Area B1 with a grade of C- or better; and completion of GE Area C1."
$ws.Range("D62").Value = "NA"
$ws.Range("E62").Value = "NA"
$ws.Range("F62").Value = "NA"
$ws.Range("G62").Value = "F, SP"
$ws.Range("C63").Value = "Junior standing; completion of GE Area A with grades of C- or better; completion of one course in GE Area B1 with a grade of C- or better; and completion of GE Area C1."
$ws.Range("D63").Value = "NA"
$ws.Range("E63").Value = "NA"
$ws.Range("F63").Value = "NA"
$ws.Range("G63").Value = "W"
$ws.Range("C64").Value = "Junior standing or English major; completion of GE Area A with grades of C- or better; completion of one course in GE Area B1 with a grade of C- or better; and completion of GE Area C1."
$ws.Range("D64").Value = "NA"
$ws.Range("E64").Value = "NA"
$ws.Range("F64").Value = "NA"
$ws.Range("G64").Value = "F, W, SP"
$ws.Range("C65").Value = "Junior standing or English major; completion of GE Area A with grades of C- or better; completion of one course in GE Area B1 with a grade of C- or better; and completion of GE Area C1."
$ws.Range("D65").Value = "NA"
$ws.Range("E65").Value = "NA"
$ws.Range("F65").Value = "NA"
$ws.Range("G65").Value = "F, W, SU"
$ws.Range("C66").Value = "Junior standing or English major; completion of GE Area A with grades of C- or better; completion of one course in GE Area B1 with a grade of C- or better; and completion of GE Area C1."
$ws.Range("D66").Value = "NA"
$ws.Range("E66").Value = "NA"
$ws.Range("F66").Value = "NA"
$ws.Range("G66").Value = "TBD"
$ws.Range("C67").Value = "Junior standing or English major; completion of GE Area A with grades of C- or better; completion of one course in GE Area B1 with a grade of C- or better; and completion of GE Area C1."
$ws.Range("D67").Value = "NA"
$ws.Range("E67").Value = "NA"
$ws.Range("F67").Value = "NA"
$ws.Range("G67").Value = "TBD"
$ws.Range("C68").Value = "Junior standing or English major; completion of GE Area A with grades of C- or better; completion of one course in GE Area B1 with a grade of C- or better; and completion of GE Area C1."
$ws.Range("D68").Value = "NA"
$ws.Range("E68").Value = "NA"
$ws.Range("F68").Value = "NA"
$ws.Range("G68").Value = "F, W, SP"
$ws.Range("C69").Value = "Junior standing or English major; completion of GE Area A with grades of C- or better; completion of one course in GE Area B1 with a grade of C- or better; and completion of GE Area C1."
$ws.Range("D69").Value = "NA"
$ws.Range("E69").Value = "NA"
$ws.Range("F69").Value = "NA"
$ws.Range("G69").Value = "F, W, SP"
$ws.Range("C70").Value = "Completion of GE Area A with grades of C- or better."
$ws.Range("D70").Value = "NA"
$ws.Range("E70").Value = "NA"
$ws.Range("F70").Value = "NA"
$ws.Range("G70").Value = "SP"
$ws.Range("C71").Value = "Junior standing and completion of GE Area A with grades of C- or better."
$ws.Range("D71").Value = "NA"
$ws.Range("E71").Value = "NA"
$ws.Range("F71").Value = "NA"
$ws.Range("G71").Value = "W, SP"
$ws.Range("C72").Value = "Completion of GE Area A with grades of C- or better."
$ws.Range("D72").Value = "NA"
$ws.Range("E72").Value = "NA"
$ws.Range("F72").Value = "NA"
$ws.Range("G72").Value = "W"
$ws.Range("C73").Value = "Consent of the department chair."
$ws.Range("D73").Value = "NA"
$ws.Range("E73").Value = "NA"
$ws.Range("F73").Value = "NA"
$ws.Range("G73").Value = "TBD"
$ws.Range("C74").Value = "Consent of instructor."
$ws.Range("D74").Value = "NA"
$ws.Range("E74").Value = "NA"
$ws.Range("F74").Value = "NA"
$ws.Range("G74").Value = "F, W, SP"
$ws.Range("C75").Value = "Junior standing; and completion of GE Area A with grades of C- or better."
$ws.Range("D75").Value = "NA"
$ws.Range("E75").Value = "NA"
$ws.Range("F75").Value = "NA"
$ws.Range("G75").Value = "F"
$ws.Range("C76").Value = "ENGL 411."
$ws.Range("D76").Value = "NA"
$ws.Range("E76").Value = "NA"
$ws.Range("F76").Value = "NA"
$ws.Range("G76").Value = "SP"
$ws.Range("C77").Value = "Junior standing; ENGL 221; and ENGL 317 or ENGL 319."
$ws.Range("D77").Value = "NA"
$ws.Range("E77").Value = "NA"
$ws.Range("F77").Value = "NA"
$ws.Range("G77").Value = "TBD"
$ws.Range("C78").Value = "Completion of GE Area A with grades of C- or better; Senior or graduate standing and admission to the teacher education program."
$ws.Range("D78").Value = "NA"
$ws.Range("E78").Value = "NA"
$ws.Range("F78").Value = "NA"
$ws.Range("G78").Value = "F"
$ws.Range("C79").Value = "Acceptance into the Single Subject Credential Program in English."
$ws.Range("D79").Value = "NA"
$ws.Range("E79").Value = "EDUC 469 or EDUC 479."
$ws.Range("F79").Value = "NA"
$ws.Range("G79").Value = "W, SP "
$ws.Range("C80").Value = "Junior standing and two of the ENGL 203, ENGL 204, ENGL 205, ENGL 303, ENGL 304, ENGL 305, ENGL 306."
$ws.Range("D80").Value = "NA"
$ws.Range("E80").Value = "NA"
$ws.Range("F80").Value = "English Major Sequence class in the relevant period."
$ws.Range("G80").Value = "F, SP "
$ws.Range("C81").Value = "Junior standing and two of the ENGL 203, ENGL 204, ENGL 205, ENGL 303, ENGL 304, ENGL 305, ENGL 306."
$ws.Range("D81").Value = "NA"
$ws.Range("E81").Value = "NA"
$ws.Range("F81").Value = "English Major Sequence class in the relevant period."
$ws.Range("G81").Value = "F, W, SP "
$ws.Range("C82").Value = "Junior standing and two of the ENGL 203, ENGL 204, ENGL 205, ENGL 303, ENGL 304, ENGL 305, ENGL 306."
$ws.Range("D82").Value = "NA"
$ws.Range("E82").Value = "NA"
$ws.Range("F82").Value = "English Major Sequence class in the relevant period."
$ws.Range("G82").Value = "F, W "
$ws.Range("C83").Value = "Senior standing; completion of the GWR; four of the ENGL 203, ENGL 204, ENGL 205, ENGL 303, ENGL 304, ENGL 305; and three of the ENGL 430, ENGL 431, ENGL 432, ENGL 439, ENGL 449, ENGL 459."
$ws.Range("D83").Value = "NA"
$ws.Range("E83").Value = "NA"
$ws.Range("F83").Value = "NA"
$ws.Range("G83").Value = "F, SP"
$ws.Range("C84").Value = "Junior standing; completion of GE Area A with grades of C- or better; and completion of GE Area C4."
$ws.Range("D84").Value = "NA"
$ws.Range("E84").Value = "NA"
$ws.Range("F84").Value = "NA"
$ws.Range("G84").Value = "W, SP"
$ws.Range("C85").Value = "Consent of instructor."
$ws.Range("D85").Value = "NA"
$ws.Range("E85").Value = "NA"
$ws.Range("F85").Value = "NA"
$ws.Range("G85").Value = "TBD"
$ws.Range("C86").Value = "ENGL 387."
$ws.Range("D86").Value = "NA"
$ws.Range("E86").Value = "NA"
$ws.Range("F86").Value = "NA"
$ws.Range("G86").Value = "F, W"
$ws.Range("C87").Value = "ENGL 388."
$ws.Range("D87").Value = "NA"
$ws.Range("E87").Value = "NA"
$ws.Range("F87").Value = "NA"
$ws.Range("G87").Value = "F, W"
$ws.Range("C88").Value = "ENGL 290 or ENGL 390."
$ws.Range("D88").Value = "NA"
$ws.Range("E88").Value = "NA"
$ws.Range("F88").Value = "NA"
$ws.Range("G88").Value = "F"
$ws.Range("C89").Value = "Two of the ENGL 290, ENGL 390, ENGL 391, ENGL 395, ENGL 495."
$ws.Range("D89").Value = "NA"
$ws.Range("E89").Value = "NA"
$ws.Range("F89").Value = "NA"
$ws.Range("G89").Value = "F"
$ws.Range("C90").Value = "ENGL 497."
$ws.Range("D90").Value = "NA"
$ws.Range("E90").Value = "NA"
$ws.Range("F90").Value = "NA"
$ws.Range("G90").Value = "W"
$ws.Range("C91").Value = "ENGL 498."
$ws.Range("D91").Value = "NA"
$ws.Range("E91").Value = "NA"
$ws.Range("F91").Value = "NA"
$ws.Range("G91").Value = "F, SP"
$ws.Range("C92").Value = "Graduate standing in English."
$ws.Range("D92").Value = "NA"
$ws.Range("E92").Value = "NA"
$ws.Range("F92").Value = "NA"
$ws.Range("G92").Value = "F"
$ws.Range("C93").Value = "Graduate standing in English."
$ws.Range("D93").Value = "NA"
$ws.Range("E93").Value = "NA"
$ws.Range("F93").Value = "NA"
$ws.Range("G93").Value = "F"
$ws.Range("C94").Value = "Graduate standing in English."
$ws.Range("D94").Value = "NA"
$ws.Range("E94").Value = "NA"
$ws.Range("F94").Value = "NA"
$ws.Range("G94").Value = "W"
$ws.Range("C95").Value = "Graduate standing in English."
$ws.Range("D95").Value = "NA"
$ws.Range("E95").Value = "NA"
$ws.Range("F95").Value = "NA"
$ws.Range("G95").Value = "SP"
$ws.Range("C96").Value = "Graduate standing in English and ENGL 505, or consent of instructor."
$ws.Range("D96").Value = "NA"
$ws.Range("E96").Value = "Teaching of ENGL 134."
$ws.Range("F96").Value = "NA"
$ws.Range("G96").Value = "F "
$ws.Range("C97").Value = "Graduate standing in English."
$ws.Range("D97").Value = "NA"
$ws.Range("E97").Value = "NA"
$ws.Range("F97").Value = "ENGL 501."
$ws.Range("G97").Value = "TBD "
$ws.Range("C98").Value = "Graduate standing in English."
$ws.Range("D98").Value = "NA"
$ws.Range("E98").Value = "NA"
$ws.Range("F98").Value = "ENGL 501."
$ws.Range("G98").Value = "SP "
$ws.Range("C99").Value = "Graduate standing in English."
$ws.Range("D99").Value = "NA"
$ws.Range("E99").Value = "NA"
$ws.Range("F99").Value = "ENGL 501."
$ws.Range("G99").Value = "W, SP "
$ws.Range("C100").Value = "Graduate standing in English."
$ws.Range("D100").Value = "NA"
$ws.Range("E100").Value = "NA"
$ws.Range("F100").Value = "ENGL 501."
$ws.Range("G100").Value = "F "
$ws.Range("C101").Value = "Graduate standing in English and 8 units of successful graduate work."
$ws.Range("D101").Value = "NA"
$ws.Range("E101").Value = "NA"
$ws.Range("F101").Value = "NA"
$ws.Range("G101").Value = "TBD"
$ws.Range("C102").Value = "Graduate standing."
$ws.Range("D102").Value = "NA"
$ws.Range("E102").Value = "NA"
$ws.Range("F102").Value = "NA"
$ws.Range("G102").Value = "TBD"
$ws.Range("C103").Value = "Graduate standing in English; and ENGL 505."
$ws.Range("D103").Value = "NA"
$ws.Range("E103").Value = "NA"
$ws.Range("F103").Value = "NA"
$ws.Range("G103").Value = "TBD"
$ws.Range("C104").Value = "Graduate standing or consent of instructor."
$ws.Range("D104").Value = "NA"
$ws.Range("E104").Value = "NA"
$ws.Range("F104").Value = "NA"
$ws.Range("G104").Value = "TBD"
$ws.Range("C105").Value = "Graduate standing in English and the permission of the graduate advisor."
$ws.Range("D105").Value = "NA"
$ws.Range("E105").Value = "NA"
$ws.Range("F105").Value = "NA"
$ws.Range("G105").Value = "TBD"
$ws.Range("C106").Value = "Consent of graduate advisor; ENGL 501; ENGL 502; ENGL 503; and ENGL 505."
$ws.Range("D106").Value = "NA"
$ws.Range("E106").Value = "NA"
$ws.Range("F106").Value = "NA"
$ws.Range("G106").Value = "F, SP"
